$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-08 Monday" "2024-04-09 Tuesday"

Replace-Text "41×43=1763" "26×50=1300"
Replace-Text "43×47=2021" "20×77=1540"
Replace-Text "85×80=6800" "14×95=1330"
Replace-Text "66×21=1386" "17×65=1105"
Replace-Text "75×97=7275" "50×29=1450"

Replace-Text "42×82=3444" "35×13=455"
Replace-Text "31×57=1767" "25×37=925"
Replace-Text "26×39=1014" "58×82=4756"
Replace-Text "33×36=1188" "49×63=3087"
Replace-Text "86×92=7912" "44×67=2948"

Replace-Text "38×26=988" "48×39=1872"
Replace-Text "17×91=1547" "33×37=1221"
Replace-Text "22×97=2134" "86×45=3870"
Replace-Text "21×52=1092" "79×63=4977"
Replace-Text "18×17=306" "98×42=4116"

Replace-Text "99×44=4356" "73×61=4453"
Replace-Text "40×51=2040" "87×21=1827"
Replace-Text "13×68=884" "60×43=2580"
Replace-Text "28×39=1092" "68×39=2652"
Replace-Text "62×40=2480" "78×79=6162"

Replace-Text "59×56=3304" "80×19=1520"
Replace-Text "91×53=4823" "94×21=1974"
Replace-Text "30×80=2400" "71×65=4615"
Replace-Text "42×70=2940" "89×77=6853"
Replace-Text "49×12=588" "73×42=3066"
